# Convert an "RRGGBB" hex string into the BGR-packed decimal value
# that Word's Font.Color property expects.
function HexToWordColor($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

$d = $word.ActiveDocument

# Builds "- <w1> <w2> ... " at the very end of the document story (a spot
# that inherits no run formatting, so the leading "- " ends up with no
# direct character formatting at all, matching the source), colors each
# word according to $wordsAndColors, then relocates that whole chunk so
# it lands right after $anchorPos (the position right after a sentence's
# closing ". "). $wordsAndColors is a list of (text, colorHexOrNull)
# pairs, one per word (including the trailing ". ").
function InsertTranslation($anchorPos, $wordsAndColors) {
    # Build the new run sequence at the safe end-of-story position
    # (one character before the final paragraph mark).
    $buildPos = $d.Content.End - 1
    $rng = $d.Range($buildPos, $buildPos)

    $rng.InsertAfter("- ")
    $rng.Collapse(0)

    foreach ($pair in $wordsAndColors) {
        $text = $pair[0]
        $colorHex = $pair[1]
        $rng.InsertAfter($text)
        if ($colorHex) {
            $rng.Font.Color = HexToWordColor $colorHex
        }
        $rng.Collapse(0)
    }

    $buildEnd = $d.Content.End - 1
    $builtLen = $buildEnd - $buildPos
    $built = $d.Range($buildPos, $buildEnd)

    # Copy the fully-formatted chunk into place at the target position...
    $target = $d.Range($anchorPos, $anchorPos)
    $target.FormattedText = $built.FormattedText

    # ...then remove the scratch copy, which has since shifted forward
    # by the length of what was just inserted ahead of it.
    $scratch = $d.Range($buildPos + $builtLen, $buildEnd + $builtLen)
    $scratch.Text = ""
}

# ---------------------------------------------------------------
# Sentence 1: "This is a text . " -> recolor, then add "- Este es
# un texto . " right after it.
# ---------------------------------------------------------------
$d.Range(0, 5).Font.Color   = HexToWordColor "C43FA4"   # This
$d.Range(5, 8).Font.Color   = HexToWordColor "4DD547"   # is
$d.Range(8, 10).Font.Color  = HexToWordColor "D1A4A7"   # a
$d.Range(10, 15).Font.Color = HexToWordColor "3808E4"   # text
$d.Range(15, 17).Font.Color = HexToWordColor "93DAA9"   # .

InsertTranslation 17 @(
    ,("Este ", "C43FA4")
    ,("es ",   "4DD547")
    ,("un ",   "D1A4A7")
    ,("texto ","3808E4")
    ,(". ",    "93DAA9")
)

# ---------------------------------------------------------------
# Sentence 2: "He eats an apple . " -> recolor, then add "- El
# come una manzana . " right after it. Locate "He " with Find
# since sentence 1's edits shifted everything after it forward.
# ---------------------------------------------------------------
$searchRng = $d.Content
$searchRng.Find.ClearFormatting()
$searchRng.Find.Execute("He ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$heStart = $searchRng.Start

$d.Range($heStart,      $heStart + 3).Font.Color  = HexToWordColor "C43FA4"   # He
$d.Range($heStart + 3,  $heStart + 8).Font.Color  = HexToWordColor "227B96"   # eats
$d.Range($heStart + 8,  $heStart + 11).Font.Color = HexToWordColor "D1A4A7"   # an
$d.Range($heStart + 11, $heStart + 17).Font.Color = HexToWordColor "3808E4"   # apple
$d.Range($heStart + 17, $heStart + 19).Font.Color = HexToWordColor "93DAA9"   # .

InsertTranslation ($heStart + 19) @(
    ,("El ",      "D1A4A7")
    ,("come ",    "227B96")
    ,("una ",     "D1A4A7")
    ,("manzana ", "3808E4")
    ,(". ",       "93DAA9")
)

$d.Save()
